$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Act (J) and Scene (K) columns with the new, more descriptive
# text (the "Dramatic Analysis" column I itself is unchanged).
$ws.Range("J2").Value = 'Act I: "The Recruitment"'
$ws.Range("K2").Value = 'Scene 1: "Ananse announces his quest to find a fool, encountering Osansa, who declines."'
$ws.Range("J3").Value = 'Act I: "The Recruitment"'
$ws.Range("K3").Value = 'Scene 1: "Ananse announces his quest to find a fool, encountering Osansa, who declines."'
$ws.Range("J4").Value = 'Act I: "The Recruitment"'
$ws.Range("K4").Value = 'Scene 1: "Ananse announces his quest to find a fool, encountering Osansa, who declines."'
$ws.Range("J5").Value = 'Act I: "The Recruitment"'
$ws.Range("K5").Value = 'Scene 1: "Ananse announces his quest to find a fool, encountering Osansa, who declines."'
$ws.Range("J6").Value = 'Act I: "The Recruitment"'
$ws.Range("K6").Value = 'Scene 2: "Ananse convinces Anene to join him despite the latter''s initial reservations."'
$ws.Range("J7").Value = 'Act I: "The Recruitment"'
$ws.Range("K7").Value = 'Scene 2: "Ananse convinces Anene to join him despite the latter''s initial reservations."'
$ws.Range("J8").Value = 'Act I: "The Recruitment"'
$ws.Range("K8").Value = 'Scene 2: "Ananse convinces Anene to join him despite the latter''s initial reservations."'
$ws.Range("J9").Value = 'Act I: "The Recruitment"'
$ws.Range("K9").Value = 'Scene 3: "Ananse and Anene work together, with Ananse attempting to trick Anene into carrying the heavier burdens."'
$ws.Range("J10").Value = 'Act I: "The Recruitment"'
$ws.Range("K10").Value = 'Scene 3: "Ananse and Anene work together, with Ananse attempting to trick Anene into carrying the heavier burdens."'
$ws.Range("J11").Value = 'Act I: "The Recruitment"'
$ws.Range("K11").Value = 'Scene 3: "Ananse and Anene work together, with Ananse attempting to trick Anene into carrying the heavier burdens."'
$ws.Range("J12").Value = 'Act I: "The Recruitment"'
$ws.Range("K12").Value = 'Scene 3: "Ananse and Anene work together, with Ananse attempting to trick Anene into carrying the heavier burdens."'
$ws.Range("J13").Value = 'Act II: "The Manipulations"'
$ws.Range("K13").Value = 'Scene 1: "The traps yield fish, but disputes arise over the division of labor and spoils."'
$ws.Range("J14").Value = 'Act II: "The Manipulations"'
$ws.Range("K14").Value = 'Scene 1: "The traps yield fish, but disputes arise over the division of labor and spoils."'
$ws.Range("J15").Value = 'Act II: "The Manipulations"'
$ws.Range("K15").Value = 'Scene 1: "The traps yield fish, but disputes arise over the division of labor and spoils."'
$ws.Range("J16").Value = 'Act II: "The Manipulations"'
$ws.Range("K16").Value = 'Scene 1: "The traps yield fish, but disputes arise over the division of labor and spoils."'
$ws.Range("J17").Value = 'Act II: "The Manipulations"'
$ws.Range("K17").Value = 'Scene 1: "The traps yield fish, but disputes arise over the division of labor and spoils."'
$ws.Range("J18").Value = 'Act II: "The Manipulations"'
$ws.Range("K18").Value = 'Scene 1: "The traps yield fish, but disputes arise over the division of labor and spoils."'
$ws.Range("J19").Value = 'Act II: "The Manipulations"'
$ws.Range("K19").Value = 'Scene 1: "The traps yield fish, but disputes arise over the division of labor and spoils."'
$ws.Range("J20").Value = 'Act II: "The Manipulations"'
$ws.Range("K20").Value = 'Scene 1: "The traps yield fish, but disputes arise over the division of labor and spoils."'
$ws.Range("J21").Value = 'Act II: "The Manipulations"'
$ws.Range("K21").Value = 'Scene 1: "The traps yield fish, but disputes arise over the division of labor and spoils."'
$ws.Range("J22").Value = 'Act II: "The Manipulations"'
$ws.Range("K22").Value = 'Scene 2: "The final escalation as Ananse pushes his schemes too far."'
$ws.Range("J23").Value = 'Act II: "The Manipulations"'
$ws.Range("K23").Value = 'Scene 2: "The final escalation as Ananse pushes his schemes too far."'
$ws.Range("J24").Value = 'Act II: "The Manipulations"'
$ws.Range("K24").Value = 'Scene 2: "The final escalation as Ananse pushes his schemes too far."'
$ws.Range("J25").Value = 'Act II: "The Manipulations"'
$ws.Range("K25").Value = 'Scene 2: "The final escalation as Ananse pushes his schemes too far."'
$ws.Range("J26").Value = 'Act II: "The Manipulations"'
$ws.Range("K26").Value = 'Scene 3: "Anene skillfully outmaneuvers Ananse, leading to a confrontation in the village."'
$ws.Range("J27").Value = 'Act II: "The Manipulations"'
$ws.Range("K27").Value = 'Scene 3: "Anene skillfully outmaneuvers Ananse, leading to a confrontation in the village."'
$ws.Range("J28").Value = 'Act II: "The Manipulations"'
$ws.Range("K28").Value = 'Scene 3: "Anene skillfully outmaneuvers Ananse, leading to a confrontation in the village."'
$ws.Range("J29").Value = 'Act II: "The Manipulations"'
$ws.Range("K29").Value = 'Scene 3: "Anene skillfully outmaneuvers Ananse, leading to a confrontation in the village."'
$ws.Range("J30").Value = 'Act II: "The Manipulations"'
$ws.Range("K30").Value = 'Scene 3: "Anene skillfully outmaneuvers Ananse, leading to a confrontation in the village."'
$ws.Range("J31").Value = 'Act III: "The Reckoning"'
$ws.Range("K31").Value = 'Scene 1: "The Chief’s intervention and judgment against Ananse."'
$ws.Range("J32").Value = 'Act III: "The Reckoning"'
$ws.Range("K32").Value = 'Scene 1: "The Chief’s intervention and judgment against Ananse."'
$ws.Range("J33").Value = 'Act III: "The Reckoning"'
$ws.Range("K33").Value = 'Scene 1: "The Chief’s intervention and judgment against Ananse."'
$ws.Range("J34").Value = 'Act III: "The Reckoning"'
$ws.Range("K34").Value = 'Scene 1: "The Chief’s intervention and judgment against Ananse."'
$ws.Range("J35").Value = 'Act III: "The Reckoning"'
$ws.Range("K35").Value = 'Scene 2: "Ananse reflects on his actions, realizing he became the fool in his quest."'
$ws.Range("J36").Value = 'Act III: "The Reckoning"'
$ws.Range("K36").Value = 'Scene 2: "Ananse reflects on his actions, realizing he became the fool in his quest."'
$ws.Range("J37").Value = 'Act III: "The Reckoning"'
$ws.Range("K37").Value = 'Scene 2: "Ananse reflects on his actions, realizing he became the fool in his quest."'
$ws.Range("J38").Value = 'Act III: "The Reckoning"'
$ws.Range("K38").Value = 'Scene 3: "Ananse resigns to his punishment, completing his tale of folly."'
$ws.Range("J39").Value = 'Act III: "The Reckoning"'
$ws.Range("K39").Value = 'Scene 3: "Ananse resigns to his punishment, completing his tale of folly."'
$ws.Range("J40").Value = 'Act III: "The Reckoning"'
$ws.Range("K40").Value = 'Scene 3: "Ananse resigns to his punishment, completing his tale of folly."'
$ws.Range("J41").Value = 'Act IV: "The Reflection"'
$ws.Range("K41").Value = 'Scene 1: "Ananse is publicly humiliated, as the Chief declares the folly of his actions."'
$ws.Range("J42").Value = 'Act IV: "The Reflection"'
$ws.Range("K42").Value = 'Scene 1: "Ananse is publicly humiliated, as the Chief declares the folly of his actions."'
$ws.Range("J43").Value = 'Act IV: "The Reflection"'
$ws.Range("K43").Value = 'Scene 1: "Ananse is publicly humiliated, as the Chief declares the folly of his actions."'
$ws.Range("J44").Value = 'Act IV: "The Reflection"'
$ws.Range("K44").Value = 'Scene 2: "Ananse resigns to his punishment, completing his tale of folly."'
$ws.Range("J45").Value = 'Act IV: "The Reflection"'
$ws.Range("K45").Value = 'Scene 2: "Ananse resigns to his punishment, completing his tale of folly."'
$ws.Range("J46").Value = 'Act IV: "The Reflection"'
$ws.Range("K46").Value = 'Scene 3: "Ananse resigns to his punishment, completing his tale of folly."'
$ws.Range("J47").Value = 'Act IV: "The Reflection"'
$ws.Range("K47").Value = 'Scene 3: "Ananse resigns to his punishment, completing his tale of folly."'
$ws.Range("J48").Value = 'Act IV: "The Reflection"'
$ws.Range("K48").Value = 'Scene 4: "The moral of the story is revealed, emphasizing justice and irony."'
$ws.Range("J49").Value = 'Act IV: "The Reflection"'
$ws.Range("K49").Value = 'Scene 4: "The moral of the story is revealed, emphasizing justice and irony."'
$ws.Range("J50").Value = 'Act IV: "The Reflection"'
$ws.Range("K50").Value = 'Scene 4: "The moral of the story is revealed, emphasizing justice and irony."'
$ws.Range("J51").Value = 'Act IV: "The Reflection"'
$ws.Range("K51").Value = 'Scene 4: "The moral of the story is revealed, emphasizing justice and irony."'
$ws.Range("J52").Value = 'Act IV: "The Reflection"'
$ws.Range("K52").Value = 'Scene 4: "The moral of the story is revealed, emphasizing justice and irony."'
$ws.Range("J53").Value = 'Act IV: "The Reflection"'
$ws.Range("K53").Value = 'Scene 5: "The story closes with Ananse''s reflection on the lessons learned, bringing justice full circle."'
$ws.Range("J54").Value = 'Act IV: "The Reflection"'
$ws.Range("K54").Value = 'Scene 5: "The story closes with Ananse''s reflection on the lessons learned, bringing justice full circle."'

# Widen column J (Act) now that its text is longer.
$ws.Columns.Item(10).ColumnWidth = 9

# Update the view: zoom out and select the newly updated Scene column.
$excel.ActiveWindow.Zoom = 80
[void]$ws.Range("K2:K54").Select()
